$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add "Bullets Per Shot" column to the table (extends Table3 from A1:I11 to A1:J11) ---
$tbl = $ws.ListObjects.Item(1)
$tbl.ListColumns.Add() | Out-Null
$ws.Range("J1").Value = "Bullets Per Shot"

# --- Column width for new column J (closest achievable to source width) ---
$ws.Columns.Item(10).ColumnWidth = 14

# --- Give new J2:J11 cells the same style (left-aligned) used by the rest of the data rows ---
$ws.Range("J2:J11").HorizontalAlignment = -4131

# --- Update "Falloff Ratio" (I) values + add "Bullets Per Shot" (J) values per weapon row ---
# Row 2: Rifle
$ws.Range("I2").Value = 0.4
$ws.Range("J2").Value = 1
# Row 3: Pistol
$ws.Range("I3").Value = 0.3
$ws.Range("J3").Value = 1
# Row 4: SMG
$ws.Range("G4").Value = 5
$ws.Range("H4").Value = 30
$ws.Range("I4").Value = 0.1
$ws.Range("J4").Value = 1
# Row 5: Revolver
$ws.Range("J5").Value = 1
# Row 6: Crossbow
$ws.Range("I6").Value = 0.1
$ws.Range("J6").Value = 1
# Row 7: SaltGun
$ws.Range("I7").Value = 0.1
$ws.Range("J7").Value = 1
# Row 8: Minigun
$ws.Range("I8").Value = 1
$ws.Range("J8").Value = 1
# Row 9: Shotgun
$ws.Range("B9").Value = 15
$ws.Range("G9").Value = 4
$ws.Range("H9").Value = 20
$ws.Range("I9").Value = 0.1
$ws.Range("J9").Value = 10
# Row 10: RocketLauncher
$ws.Range("I10").Value = 0.1
$ws.Range("J10").Value = 1

# --- Update TTK calculation formulas (rows 19-27) to account for Bullets Per Shot ---
$ws.Range("B19").Formula = '=ROUND((B2*(1/D2))*J2*$C$17,0)'
$ws.Range("G19").Formula = '=ROUND(((B2*(1/D2))*J2*$H$17)*C2*$H$16,0)'

$ws.Range("B20").Formula = '=ROUND((B3*(1/D3))*J3*$C$17,0)'
$ws.Range("B21").Formula = '=ROUND((B4*(1/D4))*J4*$C$17,0)'
$ws.Range("B22").Formula = '=ROUND((B5*(1/D5))*J5*$C$17,0)'
$ws.Range("B23").Formula = '=ROUND((B6*(1/D6))*J6*$C$17,0)'
$ws.Range("B24").Formula = '=ROUND((B7*(1/D7))*J7*$C$17,0)'
$ws.Range("B25").Formula = '=ROUND((B8*(1/D8))*J8*$C$17,0)'
$ws.Range("B26").Formula = '=ROUND((B9*(1/D9))*J9*$C$17,0)'
$ws.Range("B27").Formula = '=ROUND((B10*(1/D10))*J10*$C$17,0)'

$ws.Range("G20").Formula = '=ROUND(((B3*(1/D3))*J3*$H$17)*C3*$H$16,0)'
$ws.Range("G21").Formula = '=ROUND(((B4*(1/D4))*J4*$H$17)*C4*$H$16,0)'
$ws.Range("G22").Formula = '=ROUND(((B5*(1/D5))*J5*$H$17)*C5*$H$16,0)'
$ws.Range("G23").Formula = '=ROUND(((B6*(1/D6))*J6*$H$17)*C6*$H$16,0)'
$ws.Range("G24").Formula = '=ROUND(((B7*(1/D7))*J7*$H$17)*C7*$H$16,0)'
$ws.Range("G25").Formula = '=ROUND(((B8*(1/D8))*J8*$H$17)*C8*$H$16,0)'
$ws.Range("G26").Formula = '=ROUND(((B9*(1/D9))*J9*$H$17)*C9*$H$16,0)'
$ws.Range("G27").Formula = '=ROUND(((B10*(1/D10))*J10*$H$17)*C10*$H$16,0)'

# --- Selection moved to J10 ---
$ws.Range("J10").Select() | Out-Null
